$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in column B (Name) values in the same order the strings were
# --- first introduced, then back-fill the D (Code) / E (Code System)
# --- values, so the shared-string table grows in the same sequence as
# --- the authored workbook.

$ws.Cells.Item(104, 2).Value = "Variant Call Format result"
$ws.Cells.Item(105, 2).Value = "File Format"
$ws.Cells.Item(106, 2).Value = "Reference"
$ws.Cells.Item(106, 4).Value = "reference"
$ws.Cells.Item(105, 4).Value = "fileFormat"
$ws.Cells.Item(105, 5).Value = "VCF"
$ws.Cells.Item(107, 2).Value = "Phasing"
$ws.Cells.Item(107, 4).Value = "phasing"
$ws.Cells.Item(108, 2).Value = "Information"
$ws.Cells.Item(109, 2).Value = "Filter"
$ws.Cells.Item(110, 2).Value = "Format"
$ws.Cells.Item(108, 4).Value = "INFO"
$ws.Cells.Item(109, 4).Value = "FILTER"
$ws.Cells.Item(110, 4).Value = "FORMAT"
$ws.Cells.Item(112, 2).Value = "Number"
$ws.Cells.Item(114, 2).Value = "Description"

# --- Remaining cells (id numbers, value types, and the cells that reuse
# --- already-existing shared strings).

# Row 104 - Variant Call Format result (binary, no code / code system)
$ws.Cells.Item(104, 1).Value = 103
$ws.Cells.Item(104, 3).Value = "binary"

# Row 105 - File Format
$ws.Cells.Item(105, 1).Value = 104
$ws.Cells.Item(105, 3).Value = "short_text"

# Row 106 - Reference
$ws.Cells.Item(106, 1).Value = 105
$ws.Cells.Item(106, 3).Value = "short_text"
$ws.Cells.Item(106, 5).Value = "VCF"

# Row 107 - Phasing
$ws.Cells.Item(107, 1).Value = 106
$ws.Cells.Item(107, 3).Value = "short_text"
$ws.Cells.Item(107, 5).Value = "VCF"

# Row 108 - Information
$ws.Cells.Item(108, 1).Value = 107
$ws.Cells.Item(108, 3).Value = "binary"
$ws.Cells.Item(108, 5).Value = "VCF"

# Row 109 - Filter
$ws.Cells.Item(109, 1).Value = 108
$ws.Cells.Item(109, 3).Value = "binary"
$ws.Cells.Item(109, 5).Value = "VCF"

# Row 110 - Format
$ws.Cells.Item(110, 1).Value = 109
$ws.Cells.Item(110, 3).Value = "binary"
$ws.Cells.Item(110, 5).Value = "VCF"

# Row 111 - ID
$ws.Cells.Item(111, 1).Value = 110
$ws.Cells.Item(111, 2).Value = "ID"
$ws.Cells.Item(111, 3).Value = "short_text"
$ws.Cells.Item(111, 4).Value = "ID"
$ws.Cells.Item(111, 5).Value = "VCF"

# Row 112 - Number
$ws.Cells.Item(112, 1).Value = 111
$ws.Cells.Item(112, 3).Value = "short_text"
$ws.Cells.Item(112, 4).Value = "Number"
$ws.Cells.Item(112, 5).Value = "VCF"

# Row 113 - Type
$ws.Cells.Item(113, 1).Value = 112
$ws.Cells.Item(113, 2).Value = "Type"
$ws.Cells.Item(113, 3).Value = "short_text"
$ws.Cells.Item(113, 4).Value = "Type"
$ws.Cells.Item(113, 5).Value = "VCF"

# Row 114 - Description
$ws.Cells.Item(114, 1).Value = 113
$ws.Cells.Item(114, 3).Value = "short_text"
$ws.Cells.Item(114, 4).Value = "Description"
$ws.Cells.Item(114, 5).Value = "VCF"

# --- Extend the "INSERT INTO dbo.attributes" formula down through the
# --- newly added rows (104-114), mirroring the fill that produced
# --- G67:G103 previously.

$formula = "=CONCATENATE(""INSERT INTO dbo.attributes (id, name, value_type, code, code_system) VALUES ("", A104, "", '"", B104, ""', '"", C104, ""', "", IF(D104 = """", ""NULL"", CONCATENATE(""'"", D104, ""'"")), "", "", IF(E104 = """", ""NULL"", CONCATENATE(""'"", E104, ""'"")), "")"")"
$ws.Range("G104:G114").Formula = $formula

# --- Update the view so the sheet is scrolled near the new rows with the
# --- last of them selected, matching the authored view state.
$ws.Range("G112").Select()
